$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Price (D) and Volume(1h) (E) columns for rows with changed values ---
# Cells whose new Price value would otherwise be auto-parsed by Excel as a number
# are forced to Text format first so they keep their original string representation.
$ws.Range("D2").Value = "66.629.54"
$ws.Range("E2").Value = "  +0.76%  "
$ws.Range("D3").Value = "3.602.71"
$ws.Range("E3").Value = "  +1.29%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "609.42"
$ws.Range("E5").Value = "  +0.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.83"
$ws.Range("E6").Value = "  +3.08%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  +0.37%  "
$ws.Range("E9").Value = "  +0.27%  "
$ws.Range("E10").Value = "  +0.67%  "
$ws.Range("E11").Value = "  +0.90%  "
$ws.Range("D12").Value = "4.211.61"
$ws.Range("E12").Value = "  +1.26%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000210"
$ws.Range("E13").Value = "  +1.45%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "29.93"
$ws.Range("E14").Value = "  -0.34%  "
$ws.Range("D15").Value = "3.617.76"
$ws.Range("E15").Value = "  +1.74%  "
$ws.Range("D16").Value = "66.713.43"
$ws.Range("E16").Value = "  +0.75%  "
$ws.Range("E17").Value = "  +0.99%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.56"
$ws.Range("E18").Value = "  +1.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.37"
$ws.Range("E19").Value = "  +2.99%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.11"
$ws.Range("E20").Value = "  +1.62%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "427.63"
$ws.Range("E21").Value = "  -0.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.619"
$ws.Range("E22").Value = "  +1.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "78.92"
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("D24").Value = "3.745.79"
$ws.Range("E24").Value = "  +1.27%  "
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000124"
$ws.Range("E26").Value = "  +5.90%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.32"
$ws.Range("E27").Value = "  +4.48%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.38"
$ws.Range("E28").Value = "  +3.30%  "
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.159"
$ws.Range("E31").Value = "  +4.70%  "
$ws.Range("D32").Value = "3.598.76"
$ws.Range("E32").Value = "  +1.33%  "
$ws.Range("E33").Value = "  -0.12%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.47"
$ws.Range("E34").Value = "  -0.17%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.86"
$ws.Range("E35").Value = "  -0.26%  "
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.66"
$ws.Range("E37").Value = "  +0.97%  "
$ws.Range("E38").Value = "  -1.91%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "177.73"
$ws.Range("E39").Value = "  +4.20%  "
$ws.Range("E40").Value = "  +0.62%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.24"
$ws.Range("E41").Value = "  +1.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.900"
$ws.Range("E42").Value = "  +0.52%  "
$ws.Range("E43").Value = "  -0.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.57"
$ws.Range("E44").Value = "  +9.24%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.999"
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("E46").Value = "  -1.20%  "
$ws.Range("E49").Value = "  +0.97%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.954"
$ws.Range("E50").Value = "  +0.84%  "

# --- Rows 47 and 48 swap places (InjectiveProtocol now ranks above EnergySwap) ---
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "25.07"
$ws.Range("E47").Value = "  -3.08%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "24.33"
$ws.Range("E48").Value = "  +3.84%  "

# --- Row 51: TheGraph replaced by Maker ---
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "2.425.34"
$ws.Range("E51").Value = "  +5.26%  "
